$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in a new log entry in row 47 (Date + Start Time), matching the
# "Work Time" table pattern used by the other rows. The "End Time" (C47)
# is left blank, same as before, so the ABS(C47-B47) formula in D47 now
# evaluates to B47 instead of 0.
$ws.Range("A47").Value = 43079
$ws.Range("B47").Value = 0.069444444444444434

# Recalculate so the shared formula in D47 and the table total formula in
# D50 pick up the new values.
$excel.CalculateFull()

# Update the view: scrolled position and active selection moved further
# down the sheet (towards the newly edited rows).
$ws.Application.ActiveWindow.ScrollRow = 31
$ws.Range("I49").Select()
